$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New roster table (player, position, team) for rows 2..19
$data = @(
    @("Jalen Brunson", "PG", "New York Knicks"),
    @("Jamal Murray", "PG,SG", "Denver Nuggets"),
    @("Devin Booker", "PG,SG", "Phoenix Suns"),
    @("D'Angelo Russell", "PG", "Brooklyn Nets"),
    @("Norman Powell", "SG,SF", "LA Clippers"),
    @("Tari Eason", "SF,PF", "Houston Rockets"),
    @("Kawhi Leonard", "SG,SF,PF", "LA Clippers"),
    @("Desmond Bane", "SG,SF", "Memphis Grizzlies"),
    @("Bam Adebayo", "C", "Miami Heat"),
    @("Myles Turner", "C", "Indiana Pacers"),
    @("Trae Young", "PG", "Atlanta Hawks"),
    @("Walker Kessler", "C", "Utah Jazz"),
    @("Guerschon Yabusele", "PF,C", "Philadelphia 76ers"),
    @("Jaylen Wells", "SF", "Memphis Grizzlies"),
    @("Scoot Henderson", "PG", "Portland Trail Blazers"),
    @("Immanuel Quickley", "PG,SG", "Toronto Raptors"),
    @("Brandon Ingram", "SG,SF,PF", "New Orleans Pelicans"),
    @("LeBron James", "SF,PF", "Los Angeles Lakers")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
